# "final touches - saturday morning"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Ornstein–Uhlenbeck" is tab 2 / the active sheet

# Row 14 caption is no longer needed on this sheet - blank it out but keep formatting
$ws.Range("A14").ClearContents()

# Highlight the "winning" (smaller) value of each B/C pair in rows 19-21 in bold,
# and give the whole mini-table the 7-decimal number format used elsewhere.
$ws.Range("B19").NumberFormat = "0.0000000"
$ws.Range("C19").NumberFormat = "0.0000000"
$ws.Range("C19").Font.Bold = $true

$ws.Range("B20").NumberFormat = "0.0000000"
$ws.Range("B20").Font.Bold = $true
$ws.Range("C20").NumberFormat = "0.0000000"

$ws.Range("B21").NumberFormat = "0.0000000"
$ws.Range("B21").Font.Bold = $true
$ws.Range("C21").NumberFormat = "0.0000000"

# Fill in the missing row 22 (k=200 results) to match the rest of the table
$ws.Range("B22").Value = 0.79070413258234695
$ws.Range("B22").NumberFormat = "0.0000000"
$ws.Range("B22").Font.Bold = $true

$ws.Range("C22").Value = 0.81894339152506601
$ws.Range("C22").NumberFormat = "0.0000000"

# Leave the selection on the newly completed range
$ws.Range("C19:C22").Select() | Out-Null

# Make sure the sheet prints on A4/Letter-ish portrait paper
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
